$wb = $excel.ActiveWorkbook

# Update the existing "Database" sheet: add the information_set / process_set
# columns and reorder so the header row reads:
# A1=information_set, B1=material_entity_set, C1=named_thing_set, D1=process_set
$dbSheet = $wb.Worksheets.Item("Database")
$dbSheet.Range("A1").Value = "information_set"
$dbSheet.Range("B1").Value = "material_entity_set"
$dbSheet.Range("C1").Value = "named_thing_set"
$dbSheet.Range("D1").Value = "process_set"

# Add the new "Process" worksheet as the 4th sheet (after MaterialEntity)
$processSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$processSheet.Name = "Process"
$processSheet.Range("A1").Value = "has_inputs"
$processSheet.Range("B1").Value = "has_outputs"
$processSheet.Range("C1").Value = "id"
$processSheet.Range("D1").Value = "name"
$processSheet.Range("E1").Value = "description"

# Add the new "Information" worksheet as the 5th sheet (after Process)
$infoSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$infoSheet.Name = "Information"
$infoSheet.Range("A1").Value = "id"
$infoSheet.Range("B1").Value = "name"
$infoSheet.Range("C1").Value = "description"
